$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.85291485710751
$ws.Range("C2").Value = 7.51479035970834
$ws.Range("D2").Value = 6.017379101360685
$ws.Range("E2").Value = 12.03318672523925
$ws.Range("G2").Value = 3.687765245645171
$ws.Range("I2").Value = 28.39181915427494
$ws.Range("K2").Value = 9.52830939865915
$ws.Range("L2").Value = 10.27417356110174
$ws.Range("M2").Value = 14.11134073021696
$ws.Range("O2").Value = 29.01064532525503
$ws.Range("B3").Value = 11.65109462683541
$ws.Range("C3").Value = 7.447678406659391
$ws.Range("D3").Value = 5.90412559057639
$ws.Range("E3").Value = 12.05778784699759
$ws.Range("G3").Value = 3.689822905294766
$ws.Range("I3").Value = 28.4710514314611
$ws.Range("K3").Value = 9.380559506087367
$ws.Range("L3").Value = 10.28290663809002
$ws.Range("M3").Value = 14.08706176307511
$ws.Range("O3").Value = 29.07464507363334
$ws.Range("B4").Value = 11.52779534987723
$ws.Range("C4").Value = 7.405432920558257
$ws.Range("D4").Value = 5.835244950886624
$ws.Range("E4").Value = 12.07430536930462
$ws.Range("G4").Value = 3.691153742538503
$ws.Range("I4").Value = 28.52442035160626
$ws.Range("K4").Value = 9.290463980784143
$ws.Range("L4").Value = 10.2896803104525
$ws.Range("M4").Value = 14.07427696697764
$ws.Range("O4").Value = 29.119349148228
$ws.Range("B5").Value = 11.47777308460795
$ws.Range("C5").Value = 7.387962418062955
$ws.Range("D5").Value = 5.807383411087869
$ws.Range("E5").Value = 12.08139186687419
$ws.Range("G5").Value = 3.691713078001477
$ws.Range("I5").Value = 28.5473543697458
$ws.Range("K5").Value = 9.253952585677039
$ws.Range("L5").Value = 10.29279603694826
$ws.Range("M5").Value = 14.06960497552795
$ws.Range("O5").Value = 29.13892426729926
$ws.Range("B6").Value = 11.46948238714724
$ws.Range("C6").Value = 7.385046222866262
$ws.Range("D6").Value = 5.802770826169749
$ws.Range("E6").Value = 12.0825900536217
$ws.Range("G6").Value = 3.691806984078969
$ws.Range("I6").Value = 28.55123412851331
$ws.Range("K6").Value = 9.247903532083081
$ws.Range("L6").Value = 10.29333487838195
$ws.Range("M6").Value = 14.06886179896628
$ws.Range("O6").Value = 29.14225664848855
$ws.Range("B7").Value = 11.52711973889861
$ws.Range("C7").Value = 7.405198332233693
$ws.Range("D7").Value = 5.834868302340922
$ws.Range("E7").Value = 12.07439950055573
$ws.Range("G7").Value = 3.691161216993969
$ws.Range("I7").Value = 28.52472484828623
$ws.Range("K7").Value = 9.289970690892046
$ws.Range("L7").Value = 10.28972089066849
$ws.Range("M7").Value = 14.07421177560505
$ws.Range("O7").Value = 29.11960765011191
$ws.Range("B8").Value = 11.78323904779663
$ws.Range("C8").Value = 7.491869264532385
$ws.Range("D8").Value = 5.978218273381128
$ws.Range("E8").Value = 12.04137626316705
$ws.Range("G8").Value = 3.688460764816673
$ws.Range("I8").Value = 28.41815814773322
$ws.Range("K8").Value = 9.477263667470597
$ws.Range("L8").Value = 10.27689205733681
$ws.Range("M8").Value = 14.10253140293615
$ws.Range("O8").Value = 29.0315888436695
$ws.Range("B9").Value = 12.28747488033951
$ws.Range("C9").Value = 7.65333441360648
$ws.Range("D9").Value = 6.262695814785061
$ws.Range("E9").Value = 11.98780930510236
$ws.Range("G9").Value = 3.683697731001372
$ws.Range("I9").Value = 28.24668134155734
$ws.Range("K9").Value = 9.847463500913859
$ws.Range("L9").Value = 10.2629126494304
$ws.Range("M9").Value = 14.17471113261786
$ws.Range("O9").Value = 28.90197920262235
$ws.Range("B10").Value = 12.65532275017927
$ws.Range("C10").Value = 7.766451157741251
$ws.Range("D10").Value = 6.471339669379189
$ws.Range("E10").Value = 11.95525615548866
$ws.Range("G10").Value = 3.680519555264275
$ws.Range("I10").Value = 28.14362393081944
$ws.Range("K10").Value = 10.11855672173898
$ws.Range("L10").Value = 10.25942421365668
$ws.Range("M10").Value = 14.23760591138505
$ws.Range("O10").Value = 28.83307364015824
$ws.Range("B11").Value = 12.82132624807686
$ws.Range("C11").Value = 7.816646144200814
$ws.Range("D11").Value = 6.565692521582239
$ws.Range("E11").Value = 11.94191953679303
$ws.Range("G11").Value = 3.679142742740411
$ws.Range("I11").Value = 28.10173183401605
$ws.Range("K11").Value = 10.24114446920346
$ws.Range("L11").Value = 10.25930167117705
$ws.Range("M11").Value = 14.26829319093997
$ws.Range("O11").Value = 28.80746056776778
$ws.Range("B12").Value = 12.8839347812833
$ws.Range("C12").Value = 7.835466522965824
$ws.Range("D12").Value = 6.601302703922871
$ws.Range("E12").Value = 11.93708060232623
$ws.Range("G12").Value = 3.678631240126395
$ws.Range("I12").Value = 28.08658684824587
$ws.Range("K12").Value = 10.28741611296889
$ws.Range("L12").Value = 10.25946499383342
$ws.Range("M12").Value = 14.28020597454901
$ws.Range("O12").Value = 28.79858703598794
$ws.Range("B13").Value = 12.87046313626236
$ws.Range("C13").Value = 7.831421642528241
$ws.Range("D13").Value = 6.59363930933287
$ws.Range("E13").Value = 11.9381133585723
$ws.Range("G13").Value = 3.678740963336358
$ws.Range("I13").Value = 28.08981661179136
$ws.Range("K13").Value = 10.27745802753521
$ws.Range("L13").Value = 10.25942050605775
$ws.Range("M13").Value = 14.2776274535066
$ws.Range("O13").Value = 28.80046137418482
$ws.Range("B14").Value = 12.82648251180933
$ws.Range("C14").Value = 7.818198296247303
$ws.Range("D14").Value = 6.568624795826768
$ws.Range("E14").Value = 11.94151720074471
$ws.Range("G14").Value = 3.67910046365596
$ws.Range("I14").Value = 28.10047143765797
$ws.Range("K14").Value = 10.24495451149713
$ws.Range("L14").Value = 10.25931090972747
$ws.Range("M14").Value = 14.26926744473804
$ws.Range("O14").Value = 28.80671398465523
$ws.Range("B15").Value = 12.79950832562919
$ws.Range("C15").Value = 7.810074042517864
$ws.Range("D15").Value = 6.55328604835088
$ws.Range("E15").Value = 11.9436296673003
$ws.Range("G15").Value = 3.679321951423141
$ws.Range("I15").Value = 28.10709144889467
$ws.Range("K15").Value = 10.22502439703316
$ws.Range("L15").Value = 10.25927106493503
$ws.Range("M15").Value = 14.26418455371543
$ws.Range("O15").Value = 28.81065143787815
$ws.Range("B16").Value = 12.64444187339383
$ws.Range("C16").Value = 7.763144891859598
$ws.Range("D16").Value = 6.465158868069519
$ws.Range("E16").Value = 11.95615732814931
$ws.Range("G16").Value = 3.680610916538254
$ws.Range("I16").Value = 28.14646217784228
$ws.Range("K16").Value = 10.11052674338369
$ws.Range("L16").Value = 10.25946161827217
$ws.Range("M16").Value = 14.23564166533658
$ws.Range("O16").Value = 28.83486297946318
$ws.Range("B17").Value = 12.54892751895391
$ws.Range("C17").Value = 7.73402784442538
$ws.Range("D17").Value = 6.410924102892334
$ws.Range("E17").Value = 11.96421943000912
$ws.Range("G17").Value = 3.681419280842821
$ws.Range("I17").Value = 28.17189358035128
$ws.Range("K17").Value = 10.04006622255331
$ws.Range("L17").Value = 10.25995299812351
$ws.Range("M17").Value = 14.21865884370413
$ws.Range("O17").Value = 28.85118516383305
$ws.Range("B18").Value = 12.4938683065208
$ws.Range("C18").Value = 7.717162007073703
$ws.Range("D18").Value = 6.379679157070625
$ws.Range("E18").Value = 11.96899510632954
$ws.Range("G18").Value = 3.681890723957734
$ws.Range("I18").Value = 28.18699059355144
$ws.Range("K18").Value = 9.999472652214122
$ws.Range("L18").Value = 10.26037351544443
$ws.Range("M18").Value = 14.20908659038579
$ws.Range("O18").Value = 28.86111275358563
$ws.Range("B19").Value = 12.47520719922969
$ws.Range("C19").Value = 7.711431355513651
$ws.Range("D19").Value = 6.369092682958204
$ws.Range("E19").Value = 11.97063587693199
$ws.Range("G19").Value = 3.682051463243476
$ws.Range("I19").Value = 28.19218278978325
$ws.Range("K19").Value = 9.985718311021122
$ws.Range("L19").Value = 10.26053960238884
$ws.Range("M19").Value = 14.205879410435
$ws.Range("O19").Value = 28.86456668214426
$ws.Range("B20").Value = 12.5591082560788
$ws.Range("C20").Value = 7.737139709070622
$ws.Range("D20").Value = 6.416703003268691
$ws.Range("E20").Value = 11.96334686641756
$ws.Range("G20").Value = 3.681332557404247
$ws.Range("I20").Value = 28.16913775928593
$ws.Range("K20").Value = 10.04757407738669
$ws.Range("L20").Value = 10.2598864248276
$ws.Range("M20").Value = 14.2204464729186
$ws.Range("O20").Value = 28.84939179312794
$ws.Range("B21").Value = 12.83940804243884
$ws.Range("C21").Value = 7.82208744464635
$ws.Range("D21").Value = 6.57597569781974
$ws.Range("E21").Value = 11.94051167627855
$ws.Range("G21").Value = 3.678994602317318
$ws.Range("I21").Value = 28.09732234409107
$ws.Range("K21").Value = 10.25450598727996
$ws.Range("L21").Value = 10.25933741604684
$ws.Range("M21").Value = 14.2717151046585
$ws.Range("O21").Value = 28.8048550264049
$ws.Range("B22").Value = 13.02109662157467
$ws.Range("C22").Value = 7.876510946251816
$ws.Range("D22").Value = 6.679358348161095
$ws.Range("E22").Value = 11.92681933535089
$ws.Range("G22").Value = 3.677524099884735
$ws.Range("I22").Value = 28.05457612074458
$ws.Range("K22").Value = 10.38885632560118
$ws.Range("L22").Value = 10.26020060613936
$ws.Range("M22").Value = 14.30692197929633
$ws.Range("O22").Value = 28.78056021911172
$ws.Range("B23").Value = 12.9242832440167
$ws.Range("C23").Value = 7.847566159368147
$ws.Range("D23").Value = 6.62425840985668
$ws.Range("E23").Value = 11.93401459483323
$ws.Range("G23").Value = 3.678303690894126
$ws.Range("I23").Value = 28.07700688597869
$ws.Range("K23").Value = 10.31724672695394
$ws.Range("L23").Value = 10.25962839282678
$ws.Range("M23").Value = 14.28797805377264
$ws.Range("O23").Value = 28.79308611083156
$ws.Range("B24").Value = 12.55450599714422
$ws.Range("C24").Value = 7.735733227291589
$ws.Range("D24").Value = 6.414090560857204
$ws.Range("E24").Value = 11.96374091423681
$ws.Range("G24").Value = 3.681371744185666
$ws.Range("I24").Value = 28.17038218283504
$ws.Range("K24").Value = 10.04418003721098
$ws.Range("L24").Value = 10.25991609264575
$ws.Range("M24").Value = 14.2196376890883
$ws.Range("O24").Value = 28.85020088234366
$ws.Range("B25").Value = 12.1512548723467
$ws.Range("C25").Value = 7.610598136257184
$ws.Range("D25").Value = 6.185630404282366
$ws.Range("E25").Value = 12.00110450020864
$ws.Range("G25").Value = 3.684929601120964
$ws.Range("I25").Value = 28.28904873106398
$ws.Range("K25").Value = 9.747281663382937
$ws.Range("L25").Value = 10.26550079800168
$ws.Range("M25").Value = 14.15343048193831
$ws.Range("O25").Value = 28.93242717970074

Write-Host "Applied 240 cell updates"